$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '66.469.52'
Set-TextValue "E2" '  -0.59%  '
Set-TextValue "D3" '3.326.55'
Set-TextValue "E3" '  -0.38%  '
Set-TextValue "E4" '  -0.03%  '
Set-TextValue "D5" '586.33'
Set-TextValue "E5" '  +2.10%  '
Set-TextValue "D6" '183.52'
Set-TextValue "E6" '  +0.32%  '
Set-TextValue "D7" '0.646'
Set-TextValue "E7" '  +7.85%  '
Set-TextValue "E8" '  +0.03%  '
Set-TextValue "E9" '  -1.84%  '
Set-TextValue "D10" '6.76'
Set-TextValue "E10" '  +1.70%  '
Set-TextValue "E11" '  -0.21%  '
Set-TextValue "D12" '3.905.94'
Set-TextValue "E12" '  -0.40%  '
Set-TextValue "E13" '  -4.30%  '
Set-TextValue "D14" '66.499.37'
Set-TextValue "E14" '  -0.64%  '
Set-TextValue "D15" '26.39'
Set-TextValue "E15" '  -2.95%  '
Set-TextValue "B16" 'ShibaInu'
Set-TextValue "C16" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D16" '0.0000164'
Set-TextValue "E16" '  -2.16%  '
Set-TextValue "B17" 'WrappedEther'
Set-TextValue "C17" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D17" '3.314.70'
Set-TextValue "E17" '  -0.73%  '
Set-TextValue "D18" '431.48'
Set-TextValue "E18" '  -0.93%  '
Set-TextValue "D19" '13.33'
Set-TextValue "E19" '  -2.74%  '
Set-TextValue "E20" '  -2.85%  '
Set-TextValue "E21" '  -2.90%  '
Set-TextValue "D22" '72.32'
Set-TextValue "E22" '  -1.96%  '
Set-TextValue "E23" '  +0.14%  '
Set-TextValue "E24" '  +0.55%  '
Set-TextValue "D25" '3.454.83'
Set-TextValue "E25" '  -0.82%  '
Set-TextValue "E26" '  -0.60%  '
Set-TextValue "B27" 'Kaspa'
Set-TextValue "C27" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D27" '0.199'
Set-TextValue "E27" '  +4.08%  '
Set-TextValue "B28" 'PEPE'
Set-TextValue "C28" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D28" '0.0000115'
Set-TextValue "E28" '  -3.23%  '
Set-TextValue "E29" '  -0.56%  '
Set-TextValue "E30" '  -0.05%  '
Set-TextValue "E31" '  -0.87%  '
Set-TextValue "E32" '  -1.93%  '
Set-TextValue "D33" '1.00'
Set-TextValue "E33" '  +0.05%  '
Set-TextValue "E34" '  -2.22%  '
Set-TextValue "D35" '6.64'
Set-TextValue "E35" '  -3.16%  '
Set-TextValue "E36" '  -3.68%  '
Set-TextValue "D37" '159.95'
Set-TextValue "E37" '  -0.03%  '
Set-TextValue "E38" '  -2.80%  '
Set-TextValue "E39" '  -0.92%  '
Set-TextValue "D40" '2.894.14'
Set-TextValue "E40" '  +1.88%  '
Set-TextValue "D41" '26.72'
Set-TextValue "E41" '  -2.80%  '
Set-TextValue "E42" '  -3.29%  '
Set-TextValue "E43" '  -2.67%  '
Set-TextValue "D44" '40.34'
Set-TextValue "E44" '  +0.17%  '
Set-TextValue "E45" '  -1.31%  '
Set-TextValue "E46" '  -3.75%  '
Set-TextValue "D47" '2.33'
Set-TextValue "E47" '  -1.77%  '
Set-TextValue "D48" '23.51'
Set-TextValue "E48" '  -4.37%  '
Set-TextValue "D49" '317.92'
Set-TextValue "E49" '  -1.81%  '
Set-TextValue "E50" '  -0.51%  '
Set-TextValue "E51" '  +4.83%  '
